# Update "想去人数" (want-to-go count) values in F column across sheets
# "展览" (sheet1), "演出" (sheet2), and "全部类型" (sheet4)

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 206
$ws1.Range("F3").Value = 116
$ws1.Range("F6").Value = 5495
$ws1.Range("F7").Value = 491
$ws1.Range("F8").Value = 689
$ws1.Range("F11").Value = 78
$ws1.Range("F13").Value = 586
$ws1.Range("F17").Value = 1838
$ws1.Range("F18").Value = 1471
$ws1.Range("F19").Value = 910
$ws1.Range("F21").Value = 193
$ws1.Range("F23").Value = 548
$ws1.Range("F24").Value = 151
$ws1.Range("F27").Value = 525
$ws1.Range("F28").Value = 2886
$ws1.Range("F29").Value = 178
$ws1.Range("F32").Value = 120
$ws1.Range("F34").Value = 377
$ws1.Range("F39").Value = 290
$ws1.Range("F40").Value = 718
$ws1.Range("F43").Value = 56

# ---- Sheet "演出" ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 188
$ws2.Range("F10").Value = 12

# ---- Sheet "全部类型" ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 206
$ws4.Range("F4").Value = 116
$ws4.Range("F7").Value = 5495
$ws4.Range("F8").Value = 491
$ws4.Range("F9").Value = 689
$ws4.Range("F11").Value = 188
$ws4.Range("F16").Value = 78
$ws4.Range("F18").Value = 586
$ws4.Range("F23").Value = 1838
$ws4.Range("F24").Value = 1471
$ws4.Range("F25").Value = 910
$ws4.Range("F26").Value = 193
$ws4.Range("F29").Value = 548
$ws4.Range("F32").Value = 2886
$ws4.Range("F33").Value = 178
$ws4.Range("F36").Value = 120
$ws4.Range("F38").Value = 377
$ws4.Range("F42").Value = 290
$ws4.Range("F43").Value = 718
$ws4.Range("F45").Value = 56
$ws4.Range("F48").Value = 12
